$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 26) to the mesh-terms table so the test fixture
# can confirm duplicate mesh_id/mesh_term rows are not created.
$ws.Range("A26").Value = "C21.866.915.300.200.150"
$ws.Range("B26").Value = "Brain Concussion"
$ws.Range("C26").Value = "N"
$ws.Range("D26").Value = "N"
$ws.Range("E26").Value = "N"
$ws.Range("F26").Value = "N"
$ws.Range("G26").Value = "N"
$ws.Range("H26").Value = "N"
$ws.Range("I26").Value = "N"
$ws.Range("J26").Value = "Y"
$ws.Range("K26").Value = "N"
$ws.Range("L26").Value = "N"
$ws.Range("M26").Value = "N"
$ws.Range("N26").Value = "N"
$ws.Range("O26").Value = "N"
$ws.Range("P26").Value = "N"
$ws.Range("Q26").Value = "N"

# Columns E and O are highlighted in every data row; carry that formatting
# onto the new row the same way it appears on row 25.
$ws.Range("E25").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("O25").Copy()
$ws.Range("O26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Excel leaves the new row selected after the last entry is typed in.
$ws.Range("A26:Q26").Select()
